$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("message")
$ws.Rows("296:296").Delete()
